# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 3988
$ws1.Range("F4").Value  = 2339
$ws1.Range("F8").Value  = 16
$ws1.Range("F11").Value = 61
$ws1.Range("F13").Value = 1485
$ws1.Range("F14").Value = 263
$ws1.Range("F15").Value = 2768
$ws1.Range("F16").Value = 192

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 3988
$ws4.Range("F4").Value  = 2339
$ws4.Range("F8").Value  = 16
$ws4.Range("F12").Value = 61
$ws4.Range("F16").Value = 1485
$ws4.Range("F17").Value = 263
$ws4.Range("F18").Value = 2768
$ws4.Range("F19").Value = 192
